$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 366 (shifts the existing rows 366..487 down to 367..488,
# matching the dimension growing from A1:R487 to A1:R488).
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with the new weekly price record
# (Feria Lagunitas de Puerto Montt, Brócoli, "Primera" quality, 2023-03-03 / serial 44988).
$ws.Cells.Item(366, 1).Value = 4
$ws.Cells.Item(366, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(366, 3).Value = "Los Lagos"
$ws.Cells.Item(366, 4).Value = 44988
$ws.Cells.Item(366, 5).Value = 10
$ws.Cells.Item(366, 6).Value = 100112023
$ws.Cells.Item(366, 7).Value = "Brócoli"
$ws.Cells.Item(366, 8).Value = "Sin especificar"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 800
$ws.Cells.Item(366, 11).Value = 1600
$ws.Cells.Item(366, 12).Value = 1600
$ws.Cells.Item(366, 13).Value = 1600
$ws.Cells.Item(366, 14).Value = "$/unidad"
$ws.Cells.Item(366, 15).Value = "Región Metropolitana"
$ws.Cells.Item(366, 16).Value = 1600
$ws.Cells.Item(366, 17).Value = 1
$ws.Cells.Item(366, 18).Value = "Hortaliza"
